# Randomized subjects, body, senderName
# The "receivers" sheet holds one e-mail address per row in column A.
# This edit replaces the previous hyperlinked address list with a new,
# longer list of plain (non-hyperlinked) addresses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop every existing hyperlink on the sheet (the old rows were rendered
# as mailto: hyperlinks) and restore the affected cells to the default
# "Normal" style so the blue/underlined Hyperlink look goes away too.
$ws.Hyperlinks.Delete()
$ws.Range("A1:A7").Style = "Normal"

$addresses = @(
    "ajaygoel999@gmail.com",
    "test@chromecompete.com",
    "ajay@ajaygoel.net",
    "test@ajaygoel.org",
    "me@dropboxslideshow.com",
    "test@wordzen.com",
    "rajgoel8477@gmail.com",
    "rajanderson8477@gmail.com",
    "rajwilson8477@gmail.com",
    "briansmith8477@gmail.com",
    "oliviasmith8477@gmail.com",
    "ashsmith8477@gmail.com",
    "shellysmith8477@gmail.com",
    "ajay@madsciencekidz.com",
    "ajay2@ctopowered.com",
    "ajay@arena.tec.br",
    "ajay@daustin.co"
)

for ($i = 0; $i -lt $addresses.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $addresses[$i]
}

$ws.Range("B10").Select()
